# Auto-generated PowerPoint COM-interop edit script
$p = $ppt.ActivePresentation

# 1. Remove the trailing 11 slides (old slides 15-25), deleting from the end
#    so remaining indices are unaffected by each deletion.
for ($i = $p.Slides.Count; $i -ge 15; $i--) {
    $p.Slides.Item($i).Delete()
}

# 2. Update title + verse text on the remaining slides (3-14)
$s = $p.Slides.Item(3)
$s.Shapes.Item(2).TextFrame.TextRange.Text = 'Genesis 12:2 (KJV)'
$s.Shapes.Item(3).TextFrame.TextRange.Text = '"And I will make of thee a great nation, and I will bless thee, and make thy name great; and thou shalt be a blessing."'

$s = $p.Slides.Item(4)
$s.Shapes.Item(2).TextFrame.TextRange.Text = '2 Corinthians 5:17 (KJV)'
$s.Shapes.Item(3).TextFrame.TextRange.Text = '"Therefore if any man be in Christ, he is a new creature: old things are passed away; behold, all things are become new."'

$s = $p.Slides.Item(5)
$s.Shapes.Item(2).TextFrame.TextRange.Text = 'Revelation 21:4 (KJV)'
$s.Shapes.Item(3).TextFrame.TextRange.Text = '"And God shall wipe away all tears from their eyes; and there shall be no more death, neither sorrow, nor crying, neither shall there be any more pain: for the former things are passed away."'

$s = $p.Slides.Item(6)
$s.Shapes.Item(2).TextFrame.TextRange.Text = 'Genesis 17:1 (KJV)'
$s.Shapes.Item(3).TextFrame.TextRange.Text = '"And when Abram was ninety years old and nine, the Lord appeared to Abram, and said unto him, I am the Almighty God; walk before me, and be thou perfect."'

$s = $p.Slides.Item(7)
$s.Shapes.Item(2).TextFrame.TextRange.Text = 'Genesis 17:5 (KJV)'
$s.Shapes.Item(3).TextFrame.TextRange.Text = '"Neither shall thy name any more be called Abram, but thy name shall be Abraham; for a father of many nations have I made thee."'

$s = $p.Slides.Item(8)
$s.Shapes.Item(2).TextFrame.TextRange.Text = 'Genesis 17:15 (KJV)'
$s.Shapes.Item(3).TextFrame.TextRange.Text = '"And God said unto Abraham, As for Sarai thy wife, thou shalt not call her name Sarai, but Sarah shall her name be."'

$s = $p.Slides.Item(9)
$s.Shapes.Item(2).TextFrame.TextRange.Text = 'Romans 4:13–14 (KJV) (Part 1/2)'
$s.Shapes.Item(3).TextFrame.TextRange.Text = '"For the promise, that he should be the heir of the world, was not to Abraham, or to his seed, through the law, but through the righteousness of faith."'

$s = $p.Slides.Item(10)
$s.Shapes.Item(2).TextFrame.TextRange.Text = 'Romans 4:13–14 (KJV) (Part 2/2)'
$s.Shapes.Item(3).TextFrame.TextRange.Text = '"For if they which are of the law be heirs, faith is made void, and the promise made of none effect:"'

$s = $p.Slides.Item(11)
$s.Shapes.Item(2).TextFrame.TextRange.Text = 'Galatians 3:29 (KJV)'
$s.Shapes.Item(3).TextFrame.TextRange.Text = '"And if ye be Christ''s, then are ye Abraham''s seed, and heirs according to the promise."'

$s = $p.Slides.Item(12)
$s.Shapes.Item(2).TextFrame.TextRange.Text = 'Romans 5:19 (KJV)'
$s.Shapes.Item(3).TextFrame.TextRange.Text = '"For as by one man''s disobedience many were made sinners, so by the obedience of one shall many be made righteous."'

$s = $p.Slides.Item(13)
$s.Shapes.Item(2).TextFrame.TextRange.Text = '2 Corinthians 5:21 (KJV)'
$s.Shapes.Item(3).TextFrame.TextRange.Text = '"For he hath made him to be sin for us, who knew no sin; that we might be made the righteousness of God in him."'

$s = $p.Slides.Item(14)
$s.Shapes.Item(2).TextFrame.TextRange.Text = 'Galatians 5:4 (KJV)'
$s.Shapes.Item(3).TextFrame.TextRange.Text = '"Christ is become of no effect unto you, whosoever of you are justified by the law; ye are fallen from grace."'

